$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-08-24 Sunday" "2025-08-25 Monday"

Replace-Text "544×4=2176" "271×5=1355"
Replace-Text "984×6=5904" "924×9=8316"
Replace-Text "293×9=2637" "418×3=1254"
Replace-Text "861×3=2583" "571×6=3426"
Replace-Text "139×4=556" "155×9=1395"
Replace-Text "657×8=5256" "969×3=2907"
Replace-Text "307×9=2763" "637×2=1274"
Replace-Text "281×5=1405" "378×7=2646"
Replace-Text "845×4=3380" "839×9=7551"
Replace-Text "160×4=640" "110×7=770"
Replace-Text "500×8=4000" "586×9=5274"
Replace-Text "446×9=4014" "513×3=1539"
Replace-Text "420×4=1680" "693×6=4158"
Replace-Text "900×4=3600" "342×8=2736"
Replace-Text "154×4=616" "733×4=2932"
Replace-Text "877×6=5262" "112×9=1008"
Replace-Text "649×6=3894" "174×9=1566"
Replace-Text "796×5=3980" "237×6=1422"
Replace-Text "355×5=1775" "380×3=1140"
Replace-Text "203×9=1827" "714×9=6426"
Replace-Text "602×9=5418" "133×9=1197"
Replace-Text "122×2=244" "687×5=3435"
Replace-Text "241×5=1205" "554×2=1108"
Replace-Text "738×9=6642" "724×5=3620"

Write-Output "Done"
